# Adds all IG authors as contact
# The "Metadata" sheet has two existing "Contact" rows (row 10 and 11).
# This change duplicates that Contact row two more times (rows 12 and 13),
# pushing the remaining metadata rows (Jurisdiction ... Derivation) down by
# two rows, and refreshes the "Date" value to reflect the re-export time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert two blank rows right after the existing Contact rows (before the
# Jurisdiction row), shifting everything below down.
$ws.Range("A12:A13").EntireRow.Insert()

# Copy the formatting (fill/border/font/alignment) from the last Contact
# row onto the two newly inserted rows so they match the rest of the table.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)

# Populate the two new Contact rows with the same Property/Value pair as
# the existing Contact rows.
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"

# The Date property reflects the moment the IG was (re-)exported.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"
